$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1029.85
$ws.Range("I98").Value = 1095.4445
$ws.Range("J98").Value = 439.5
$ws.Range("K98").Value = 1095.4445
$ws.Range("L98").Value = 439.5
$ws.Range("M98").Value = 402.5554999999999
$ws.Range("N98").Value = -3435.5
$ws.Range("H122").Value = 1029.85
$ws.Range("I122").Value = 1095.4445
$ws.Range("J122").Value = 439.5
$ws.Range("K122").Value = 3286.3335
$ws.Range("L122").Value = 1318.5
$ws.Range("M122").Value = -836.3335000000002
$ws.Range("N122").Value = -6218.5
$ws.Range("H132").Value = 3943.442
$ws.Range("I132").Value = 3282.4146
$ws.Range("K132").Value = 9847.2438
$ws.Range("M132").Value = -7317.2438
$ws.Range("H141").Value = 6746
$ws.Range("I141").Value = 6746
$ws.Range("K141").Value = 20238
$ws.Range("M141").Value = -15058

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2943.1738
$ws.Range("I61").Value = 2849.7273
$ws.Range("J61").Value = 4999
$ws.Range("K61").Value = 2849.7273
$ws.Range("L61").Value = 4999
$ws.Range("M61").Value = -2637.7273
$ws.Range("N61").Value = -5423
$ws.Range("H88").Value = 1515.4286
$ws.Range("I88").Value = 1028.6666
$ws.Range("J88").Value = 1880.5
$ws.Range("K88").Value = 1028.6666
$ws.Range("L88").Value = 1880.5
$ws.Range("M88").Value = -622.6666
$ws.Range("N88").Value = -2692.5
$ws.Range("H91").Value = 1515.4286
$ws.Range("I91").Value = 1028.6666
$ws.Range("J91").Value = 1880.5
$ws.Range("K91").Value = 1028.6666
$ws.Range("L91").Value = 1880.5
$ws.Range("M91").Value = 375.3334
$ws.Range("N91").Value = -4688.5
$ws.Range("H97").Value = 8844.941000000001
$ws.Range("I97").Value = 12126.5
$ws.Range("J97").Value = 4157
$ws.Range("K97").Value = 12126.5
$ws.Range("L97").Value = 4157
$ws.Range("M97").Value = -11630.5
$ws.Range("N97").Value = -5149
$ws.Range("H132").Value = 39048.965
$ws.Range("I132").Value = 41541.76
$ws.Range("K132").Value = 124625.28
$ws.Range("M132").Value = -122095.28
$ws.Range("H136").Value = 2943.1738
$ws.Range("I136").Value = 2849.7273
$ws.Range("J136").Value = 4999
$ws.Range("K136").Value = 8549.1819
$ws.Range("L136").Value = 14997
$ws.Range("M136").Value = -5999.1819
$ws.Range("N136").Value = -20097

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2100
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("H89").Value = 2100
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("H94").Value = 5001.8335
$ws.Range("I94").Value = 5533.6665
$ws.Range("K94").Value = 5533.6665
$ws.Range("M94").Value = -5082.6665
$ws.Range("H99").Value = 35711.645
$ws.Range("I99").Value = 48145.453
$ws.Range("K99").Value = 48145.453
$ws.Range("M99").Value = -46647.453
$ws.Range("H105").Value = 3400.3901
$ws.Range("I105").Value = 3321.5527
$ws.Range("K105").Value = 3321.5527
$ws.Range("M105").Value = -1574.5527
$ws.Range("H134").Value = 1512.6786
$ws.Range("I134").Value = 1544.2593
$ws.Range("K134").Value = 4632.7779
$ws.Range("M134").Value = -2097.7779
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 883.4666999999999
$ws.Range("I16").Value = 827.0769
$ws.Range("J16").Value = 1250
$ws.Range("K16").Value = 827.0769
$ws.Range("L16").Value = 1250
$ws.Range("M16").Value = -540.0769
$ws.Range("N16").Value = -1824
$ws.Range("H31").Value = 5715.75
$ws.Range("I31").Value = 3442.3333
$ws.Range("J31").Value = 7079.8
$ws.Range("K31").Value = 3442.3333
$ws.Range("L31").Value = 7079.8
$ws.Range("M31").Value = -3147.3333
$ws.Range("N31").Value = -7669.8
$ws.Range("H34").Value = 5715.75
$ws.Range("I34").Value = 3442.3333
$ws.Range("J34").Value = 7079.8
$ws.Range("K34").Value = 3442.3333
$ws.Range("L34").Value = 7079.8
$ws.Range("M34").Value = -3240.3333
$ws.Range("N34").Value = -7483.8
$ws.Range("H113").Value = 883.4666999999999
$ws.Range("I113").Value = 827.0769
$ws.Range("J113").Value = 1250
$ws.Range("K113").Value = 827.0769
$ws.Range("L113").Value = 1250
$ws.Range("M113").Value = 1342.9231
$ws.Range("N113").Value = -5590
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3599.3333
$ws.Range("J34").Value = 3879.2
$ws.Range("L34").Value = 11637.6
$ws.Range("N34").Value = -11805.6
$ws.Range("H37").Value = 62028.43
$ws.Range("J37").Value = 62028.43
$ws.Range("L37").Value = 186085.29
$ws.Range("N37").Value = -186309.29
$ws.Range("H55").Value = 5234.1665
$ws.Range("J55").Value = 5281
$ws.Range("L55").Value = 15843
$ws.Range("N55").Value = -16197
$ws.Range("H126").Value = 3210
$ws.Range("I126").Value = 3210
$ws.Range("K126").Value = 9630
$ws.Range("M126").Value = -4690

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 501750
$ws.Range("J20").Value = 3500
$ws.Range("L20").Value = 3500
$ws.Range("N20").Value = -3990
$ws.Range("H113").Value = 4079.2
$ws.Range("I113").Value = 3499
$ws.Range("K113").Value = 3499
$ws.Range("M113").Value = -1329
$ws.Range("H126").Value = 5155.476
$ws.Range("I126").Value = 4216.625
$ws.Range("J126").Value = 8159.8
$ws.Range("K126").Value = 12649.875
$ws.Range("L126").Value = 24479.4
$ws.Range("M126").Value = -10179.875
$ws.Range("N126").Value = -29419.4
$ws.Range("H132").Value = 66619.31
$ws.Range("I132").Value = 73993.71000000001
$ws.Range("J132").Value = 14998.5
$ws.Range("K132").Value = 221981.13
$ws.Range("L132").Value = 44995.5
$ws.Range("M132").Value = -219451.13
$ws.Range("N132").Value = -50055.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1338.0834
$ws.Range("I93").Value = 679.125
$ws.Range("J93").Value = 2656
$ws.Range("K93").Value = 679.125
$ws.Range("L93").Value = 2656
$ws.Range("M93").Value = 568.875
$ws.Range("N93").Value = -5152
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H136").Value = 3784.9614
$ws.Range("I136").Value = 2897.4443
$ws.Range("K136").Value = 8692.332900000001
$ws.Range("M136").Value = -6142.332900000001

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 68999
$ws.Range("J46").Value = 68999
$ws.Range("L46").Value = 68999
$ws.Range("N46").Value = -69461
$ws.Range("H96").Value = 2336
$ws.Range("J96").Value = 999
$ws.Range("L96").Value = 999
$ws.Range("N96").Value = -3745
$ws.Range("H100").Value = 891
$ws.Range("I100").Value = 754.875
$ws.Range("J100").Value = 1163.25
$ws.Range("K100").Value = 1509.75
$ws.Range("L100").Value = 2326.5
$ws.Range("M100").Value = -968.75
$ws.Range("N100").Value = -3408.5
$ws.Range("H122").Value = 637
$ws.Range("I122").Value = 688.46155
$ws.Range("J122").Value = 302.5
$ws.Range("K122").Value = 2065.38465
$ws.Range("L122").Value = 907.5
$ws.Range("M122").Value = 384.61535
$ws.Range("N122").Value = -5807.5
$ws.Range("H132").Value = 57170.42
$ws.Range("I132").Value = 63995.37
$ws.Range("J132").Value = 4845.8335
$ws.Range("K132").Value = 191986.11
$ws.Range("L132").Value = 14537.5005
$ws.Range("M132").Value = -189456.11
$ws.Range("N132").Value = -19597.5005
$ws.Range("H134").Value = 68999
$ws.Range("J134").Value = 68999
$ws.Range("L134").Value = 206997
$ws.Range("N134").Value = -212067
$ws.Range("H136").Value = 2766.923
$ws.Range("I136").Value = 2531.9644
$ws.Range("J136").Value = 3365
$ws.Range("K136").Value = 7595.8932
$ws.Range("L136").Value = 10095
$ws.Range("M136").Value = -5045.8932
$ws.Range("N136").Value = -15195

Write-Host "Applied all profit-sheet updates"